$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.352.26'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.667.64'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("E4").Value = '  +0.98%  '
$ws.Range("D5").Value = '''219.55'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '''0.5345'
$ws.Range("E6").Value = '  +1.29%  '
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("D8").Value = '''0.2664'
$ws.Range("E8").Value = '  +2.57%  '
$ws.Range("D9").Value = '''0.06393'
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("D10").Value = '''20.88'
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("D11").Value = '''0.07840'
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Value = '''4.558'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = '1.651.80'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '1.895.82'
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").Value = '''0.5545'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '0.0₅8187'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '''65.99'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '26.377.49'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = '''4.684'
$ws.Range("E20").Value = '  +2.45%  '
$ws.Range("D21").Value = '''195.41'
$ws.Range("E21").Value = '  +2.62%  '
$ws.Range("D22").Value = '''10.27'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = '''1.012'
$ws.Range("E24").Value = '  +0.88%  '
$ws.Range("D25").Value = '''146.36'
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").Value = '''7.229'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '''1.502'
$ws.Range("E29").Value = '  +4.12%  '
$ws.Range("D30").Value = '''0.05858'
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("D31").Value = '''1.283'
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("D32").Value = '''3.588'
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("D34").Value = '''1.614'
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("D35").Value = '''0.9717'
$ws.Range("E35").Value = '  +3.11%  '
$ws.Range("D36").Value = '''2.836'
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = '''2.422'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").Value = '''0.5828'
$ws.Range("E38").Value = '  +1.23%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '1.076.38'
$ws.Range("E40").Value = '  +4.76%  '
$ws.Range("D41").Value = '''0.8645'
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("D42").Value = '''5.857'
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("D44").Value = '''104.28'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '''1.012'
$ws.Range("E47").Value = '  +1.03%  '
$ws.Range("D48").Value = '''0.4394'
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").Value = '''8.063'
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("D50").Value = '0.0₈103'
$ws.Range("E50").Value = '  -8.17%  '
$ws.Range("D51").Value = '''0.05168'
$ws.Range("E51").Value = '  +0.57%  '
